$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newValues = @{
    2  = "320018813081"
    3  = "320018813092"
    4  = "320018813129"
    5  = "320018813140"
    6  = "320018813184"
    7  = "320018813200"
    8  = "320018813232"
    9  = "320018813254"
    10 = "320018813287"
    11 = "320018813302"
    12 = "320018813346"
    13 = "320018813449"
    14 = "320018813471"
    15 = "320018813493"
    16 = "320018813520"
    17 = "320018813541"
    18 = "320018813585"
    19 = "320018813600"
    20 = "320018813850"
    21 = "320018813872"
    22 = "320018813909"
    23 = "320018813910"
    24 = "320018813920"
    25 = "320018813931"
    26 = "320018813942"
}

foreach ($row in $newValues.Keys) {
    $cell = $ws.Range("P$row")
    $cell.NumberFormat = "@"
    $cell.Value = $newValues[$row]
}
